$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: CaseNo, Name, Charge, Statute, Level, Plea, Finding, (H blank), Days
$data = @(
    @("21CRB01291", "Hemmeter", "PERMISSION REQ'D TO USE LICENSED DOCK", "1501:46-12-04", "MM", "No Contest", "Guilty", $null, "25"),
    @("19CRB01525", "Pelanda", "ASSAULT - M1", "2903.13(A)", "M1", "Guilty", "Guilty", $null, "25"),
    @("19CRB01525", "Pelanda", "AGGRAVATED MENACING", "2903.21", "M1", "Guilty", "Guilty", $null, "15"),
    @("19CRB01525", "Pelanda", "DISORDERLY CONDUCT", "2917.11A1", "MM", "Guilty", "Guilty", $null, "0"),
    @("03TRD13368", "Hemmeter", "SPEED REDUCED ZONE 3RD OR MORE", "4511.21C***", "M3", "Guilty", "Guilty", $null, "0")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $cell = $ws.Cells.Item($row, $c + 1)
            # Force text storage only for values that look numeric (case
            # numbers, statute codes, day counts) so they aren't silently
            # reinterpreted as numbers on write.
            if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $val
        }
    }
}

$ws.Range("B2").Select()
